# fead: added alphabetical sorting in populating the excel file
#
# Re-sorts the attendance list alphabetically by name:
#   roche (row2, no time)      -> cyrus (row2, no time)
#   kiefer (row3, PRESENT)     -> kiefer (row3, PRESENT, updated time)
#   cyrus (row4, PRESENT)      -> roche (row4, no time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: roche -> cyrus (stays without attendance time)
$ws.Range("A2").Value = "cyrus"

# Row 3: kiefer keeps PRESENT, time updates
$ws.Range("C3").Value = "21:09:27"

# Row 4: cyrus -> roche, and clear its PRESENT/time (no longer marked present)
$ws.Range("A4").Value = "roche"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
